$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value2 = 5064.3335
$ws.Range("I31").Value2 = 1096.75
$ws.Range("K31").Value2 = 3290.25
$ws.Range("M31").Value2 = -3060.25

# Row 74
$ws.Range("H74").Value2 = 4943.4614
$ws.Range("I74").Value2 = 4943.4614
$ws.Range("K74").Value2 = 4943.4614
$ws.Range("M74").Value2 = -4007.4614

# Row 77
$ws.Range("H77").Value2 = 4943.4614
$ws.Range("I77").Value2 = 4943.4614
$ws.Range("K77").Value2 = 24717.307
$ws.Range("M77").Value2 = -20037.307

# Row 116
$ws.Range("H116").Value2 = 12584.228
$ws.Range("I116").Value2 = 4494.625
$ws.Range("J116").Value2 = 17206.857
$ws.Range("K116").Value2 = 4494.625
$ws.Range("L116").Value2 = 17206.857
$ws.Range("M116").Value2 = -1052.625
$ws.Range("N116").Value2 = -24090.857


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 7848.136
$ws.Range("I32").Value2 = 7584.2617
$ws.Range("J32").Value2 = 25000
$ws.Range("K32").Value2 = 7584.2617
$ws.Range("L32").Value2 = 25000
$ws.Range("M32").Value2 = -7297.2617
$ws.Range("N32").Value2 = -25574

# Row 74
$ws.Range("H74").Value2 = 1871.174
$ws.Range("I74").Value2 = 1716.85
$ws.Range("K74").Value2 = 1716.85
$ws.Range("M74").Value2 = -842.8499999999999

# Row 77
$ws.Range("H77").Value2 = 1871.174
$ws.Range("I77").Value2 = 1716.85
$ws.Range("K77").Value2 = 8584.25
$ws.Range("M77").Value2 = -4216.25

# Row 102
$ws.Range("H102").Value2 = 2407.0557
$ws.Range("I102").Value2 = 2208.4666
$ws.Range("J102").Value2 = 3400
$ws.Range("K102").Value2 = 2208.4666
$ws.Range("L102").Value2 = 3400
$ws.Range("M102").Value2 = -586.4666000000002
$ws.Range("N102").Value2 = -6644

# Row 103
$ws.Range("H103").Value2 = 92446
$ws.Range("J103").Value2 = 92446
$ws.Range("L103").Value2 = 92446
$ws.Range("N103").Value2 = -94790


$ws = $wb.Worksheets.Item("BSM")
# Row 88
$ws.Range("H88").Value2 = 22037.285
$ws.Range("J88").Value2 = 22037.285
$ws.Range("L88").Value2 = 22037.285
$ws.Range("N88").Value2 = -22849.285

# Row 91
$ws.Range("H91").Value2 = 22037.285
$ws.Range("J91").Value2 = 22037.285
$ws.Range("L91").Value2 = 22037.285
$ws.Range("N91").Value2 = -24845.285

# Row 99
$ws.Range("H99").Value2 = 3212.5
$ws.Range("I99").Value2 = 2950
$ws.Range("J99").Value2 = 3300
$ws.Range("K99").Value2 = 2950
$ws.Range("L99").Value2 = 3300
$ws.Range("M99").Value2 = -1452
$ws.Range("N99").Value2 = -6296


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 47623120
$ws.Range("J31").Value2 = 4622
$ws.Range("L31").Value2 = 4622
$ws.Range("N31").Value2 = -5212

# Row 34
$ws.Range("H34").Value2 = 47623120
$ws.Range("J34").Value2 = 4622
$ws.Range("L34").Value2 = 4622
$ws.Range("N34").Value2 = -5026

# Row 132
$ws.Range("H132").Value2 = 2761.5264
$ws.Range("I132").Value2 = 2797.6667
$ws.Range("K132").Value2 = 8393.000100000001
$ws.Range("M132").Value2 = -5863.000100000001


$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value2 = 14478
$ws.Range("I13").Value2 = 5050.5
$ws.Range("J13").Value2 = 33333
$ws.Range("K13").Value2 = 15151.5
$ws.Range("L13").Value2 = 99999
$ws.Range("M13").Value2 = -14983.5
$ws.Range("N13").Value2 = -100335

# Row 55
$ws.Range("H55").Value2 = 6770.6875
$ws.Range("I55").Value2 = 2833.3333
$ws.Range("J55").Value2 = 18582.75
$ws.Range("K55").Value2 = 8499.999899999999
$ws.Range("L55").Value2 = 55748.25
$ws.Range("M55").Value2 = -8322.999899999999
$ws.Range("N55").Value2 = -56102.25

# Row 131
$ws.Range("H131").Value2 = 4443.0386
$ws.Range("I131").Value2 = 2802.5
$ws.Range("J131").Value2 = 5172.1665
$ws.Range("K131").Value2 = 8407.5
$ws.Range("L131").Value2 = 15516.4995
$ws.Range("M131").Value2 = -3367.5
$ws.Range("N131").Value2 = -25596.4995

# Row 139
$ws.Range("H139").Value2 = 5050.3125
$ws.Range("I139").Value2 = 2606.3635
$ws.Range("J139").Value2 = 10427
$ws.Range("K139").Value2 = 7819.0905
$ws.Range("L139").Value2 = 31281
$ws.Range("M139").Value2 = -2679.0905
$ws.Range("N139").Value2 = -41561

# Row 140
$ws.Range("H140").Value2 = 2973.1316
$ws.Range("I140").Value2 = 1098.4193
$ws.Range("K140").Value2 = 3295.2579
$ws.Range("M140").Value2 = 1884.7421


$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value2 = 1275795
$ws.Range("I122").Value2 = 2067064
$ws.Range("J122").Value2 = 9764.5
$ws.Range("K122").Value2 = 6201192
$ws.Range("L122").Value2 = 29293.5
$ws.Range("M122").Value2 = -6198742
$ws.Range("N122").Value2 = -34193.5


$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value2 = 7320.41
$ws.Range("I7").Value2 = 7082.84
$ws.Range("K7").Value2 = 7082.84
$ws.Range("M7").Value2 = -6970.84

# Row 22
$ws.Range("H22").Value2 = 18876500
$ws.Range("I22").Value2 = 44038332
$ws.Range("J22").Value2 = 5124.75
$ws.Range("K22").Value2 = 44038332
$ws.Range("L22").Value2 = 5124.75
$ws.Range("M22").Value2 = -44038037
$ws.Range("N22").Value2 = -5714.75

# Row 27
$ws.Range("H27").Value2 = 18876500
$ws.Range("I27").Value2 = 44038332
$ws.Range("J27").Value2 = 5124.75
$ws.Range("K27").Value2 = 44038332
$ws.Range("L27").Value2 = 5124.75
$ws.Range("M27").Value2 = -44038225
$ws.Range("N27").Value2 = -5338.75

# Row 55
$ws.Range("H55").Value2 = 1138.6
$ws.Range("I55").Value2 = 616.1111
$ws.Range("J55").Value2 = 1566.091
$ws.Range("K55").Value2 = 616.1111
$ws.Range("L55").Value2 = 1566.091
$ws.Range("M55").Value2 = -443.1111
$ws.Range("N55").Value2 = -1912.091

# Row 100
$ws.Range("H100").Value2 = 20858472
$ws.Range("I100").Value2 = 3244.5715
$ws.Range("K100").Value2 = 3244.5715
$ws.Range("M100").Value2 = -2703.5715

# Row 126
$ws.Range("H126").Value2 = 7320.41
$ws.Range("I126").Value2 = 7082.84
$ws.Range("K126").Value2 = 21248.52
$ws.Range("M126").Value2 = -18778.52

# Row 128
$ws.Range("H128").Value2 = 68329.664
$ws.Range("J128").Value2 = 68329.664
$ws.Range("L128").Value2 = 68329.664
$ws.Range("N128").Value2 = -78289.664

# Row 132
$ws.Range("H132").Value2 = 4677.76
$ws.Range("I132").Value2 = 2314.5833
$ws.Range("K132").Value2 = 6943.749899999999
$ws.Range("M132").Value2 = -4413.749899999999

# Row 136
$ws.Range("H136").Value2 = 4702.5
$ws.Range("I136").Value2 = 2315.2
$ws.Range("K136").Value2 = 6945.599999999999
$ws.Range("M136").Value2 = -4395.599999999999


$ws = $wb.Worksheets.Item("WVR")
# Row 38
$ws.Range("H38").Value2 = 90000
$ws.Range("I38").Value2 = 0
$ws.Range("J38").Value2 = 90000
$ws.Range("K38").Value2 = 0
$ws.Range("L38").Value2 = 90000
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value2 = -90946

# Row 122
$ws.Range("H122").Value2 = 1428.5172
$ws.Range("I122").Value2 = 1122.4286
$ws.Range("K122").Value2 = 3367.2858
$ws.Range("M122").Value2 = -917.2857999999997

# Row 132
$ws.Range("H132").Value2 = 1431395.1
$ws.Range("I132").Value2 = 3104.5557
$ws.Range("J132").Value2 = 4002318
$ws.Range("K132").Value2 = 9313.667099999999
$ws.Range("L132").Value2 = 12006954
$ws.Range("M132").Value2 = -6783.667099999999
$ws.Range("N132").Value2 = -12012014

# Row 141
$ws.Range("H141").Value2 = 93882.73
$ws.Range("J141").Value2 = 93882.73
$ws.Range("L141").Value2 = 93882.73
$ws.Range("N141").Value2 = -104242.73

